$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 52
$ws.Range("A3").Value = 21
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 1
$ws.Range("A5").Value = 32
$ws.Range("A7").Value = 47
$ws.Range("C7").Value = "Biltwell Event Center"
$ws.Range("D7").Value = 4.7
$ws.Range("E7").Value = 373
$ws.Range("A8").Value = 48
$ws.Range("C8").Value = "CC's Exotic Pets"
$ws.Range("D8").Value = 4.1
$ws.Range("E8").Value = 42
$ws.Range("A9").Value = 7
$ws.Range("C9").Value = "Cabinet and Stone Expo"
$ws.Range("D9").Value = 4.6
$ws.Range("E9").Value = 31
$ws.Range("A10").Value = 19
$ws.Range("C10").Value = "Circle City Auto Parts"
$ws.Range("D10").Value = 4.9
$ws.Range("E10").Value = 18
$ws.Range("A11").Value = 45
$ws.Range("C11").Value = "Circle City Coatings"
$ws.Range("D11").Value = 4.6
$ws.Range("E11").Value = 19
$ws.Range("A12").Value = 34
$ws.Range("C12").Value = "Convention center"
$ws.Range("E12").Value = 473
$ws.Range("A13").Value = 29
$ws.Range("C13").Value = "Core & Main"
$ws.Range("D13").Value = 4.3
$ws.Range("E13").Value = 3
$ws.Range("A14").Value = 0
$ws.Range("C14").Value = "Creation Evidence Expo"
$ws.Range("D14").Value = 4.8
$ws.Range("E14").Value = 5
$ws.Range("A15").Value = 56
$ws.Range("C15").Value = "Curtain Call Dance Center"
$ws.Range("D15").Value = 4.4
$ws.Range("E15").Value = 9
$ws.Range("A16").Value = 5
$ws.Range("C16").Value = "Curvature Expo"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("A17").Value = 16
$ws.Range("C17").Value = "Don Mitchell Pro Shop"
$ws.Range("D17").Value = 3.3
$ws.Range("E17").Value = 7
$ws.Range("A18").Value = 44
$ws.Range("C18").Value = "Don Mitchell Pro Shops"
$ws.Range("D18").Value = 4.4
$ws.Range("E18").Value = 16
$ws.Range("A19").Value = 54
$ws.Range("C19").Value = "Elegant Stylez"
$ws.Range("D19").Value = 3.1
$ws.Range("E19").Value = 3696
$ws.Range("A21").Value = 31
$ws.Range("A22").Value = 50
$ws.Range("A23").Value = 39
$ws.Range("A24").Value = 57
$ws.Range("C24").Value = "Great Day Tattoo"
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 18
$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "Hamilton County Fairgrounds"
$ws.Range("D25").Value = 4.6
$ws.Range("E25").Value = 68
$ws.Range("A26").Value = 25
$ws.Range("C26").Value = "Harvest Pavillion"
$ws.Range("E26").Value = 18
$ws.Range("A27").Value = 58
$ws.Range("C27").Value = "INKSTINCT TATTOO"
$ws.Range("D27").Value = 4.8
$ws.Range("E27").Value = 116
$ws.Range("A28").Value = 4
$ws.Range("C28").Value = "Indiana Black Expo Inc"
$ws.Range("D28").Value = 4.3
$ws.Range("E28").Value = 39
$ws.Range("A29").Value = 28
$ws.Range("C29").Value = "Indiana Convention Center"
$ws.Range("D29").Value = 4.5
$ws.Range("E29").Value = 528
$ws.Range("A30").Value = 12
$ws.Range("C30").Value = "Indiana Fishing Expo"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("A31").Value = 27
$ws.Range("C31").Value = "Indiana Flower & Patio Show"
$ws.Range("D31").Value = 4.3
$ws.Range("E31").Value = 60
$ws.Range("A32").Value = 6
$ws.Range("C32").Value = "Indiana Latino Expo"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("A33").Value = 22
$ws.Range("C33").Value = "Indiana State Fairgrounds & Event Center"
$ws.Range("D33").Value = 4.4
$ws.Range("E33").Value = 1344
$ws.Range("A34").Value = 43
$ws.Range("C34").Value = "Indiana State Numismatic Association"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("A35").Value = 35
$ws.Range("C35").Value = "Indianapolis Auto Show"
$ws.Range("D35").Value = 3.1
$ws.Range("E35").Value = 51
$ws.Range("A36").Value = 15
$ws.Range("C36").Value = "Indianapolis Chapter of Indiana Black Expo, Inc."
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("A37").Value = 30
$ws.Range("C37").Value = "Indianapolis Competition Products"
$ws.Range("A38").Value = 24
$ws.Range("C38").Value = "Indianapolis Motor Speedway"
$ws.Range("D38").Value = 4.8
$ws.Range("E38").Value = 11013
$ws.Range("A39").Value = 14
$ws.Range("C39").Value = "Indy Air Expo"
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("A40").Value = 11
$ws.Range("C40").Value = "Laser Storm"
$ws.Range("D40").Value = 4
$ws.Range("E40").Value = 51
$ws.Range("A41").Value = 26
$ws.Range("C41").Value = "Marketplace Events - Indianapolis Office"
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 1
$ws.Range("A42").Value = 8
$ws.Range("C42").Value = "Nail Expo"
$ws.Range("D42").Value = 3.2
$ws.Range("E42").Value = 62
$ws.Range("A43").Value = 9
$ws.Range("C43").Value = "National Expo, Inc"
$ws.Range("D43").Value = 5
$ws.Range("A44").Value = 13
$ws.Range("C44").Value = "Off Road Expo"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("A45").Value = 59
$ws.Range("C45").Value = "Pan Am Tower"
$ws.Range("D45").Value = 4.3
$ws.Range("E45").Value = 28
$ws.Range("A46").Value = 40
$ws.Range("C46").Value = "Premier Surface"
$ws.Range("A47").Value = 51
$ws.Range("C47").Value = "Purdue Extension / Horticulture Building"
$ws.Range("D47").Value = 4.4
$ws.Range("E47").Value = 36
$ws.Range("A48").Value = 1
$ws.Range("C48").Value = "Royal Pin Expo"
$ws.Range("D48").Value = 4.2
$ws.Range("E48").Value = 1025
$ws.Range("A49").Value = 37
$ws.Range("C49").Value = "Royal Pin Western"
$ws.Range("D49").Value = 4.4
$ws.Range("E49").Value = 1182
$ws.Range("A50").Value = 46
$ws.Range("C50").Value = "Royal Pin Woodland"
$ws.Range("D50").Value = 4.3
$ws.Range("E50").Value = 1757
$ws.Range("A51").Value = 17
$ws.Range("C51").Value = "Samps Hack Shack Brownsburg"
$ws.Range("D51").Value = 4.9
$ws.Range("E51").Value = 22
$ws.Range("A52").Value = 38
$ws.Range("C52").Value = "Samps Hack Shack Plainfield"
$ws.Range("D52").Value = 5
$ws.Range("E52").Value = 9
$ws.Range("A53").Value = 20
$ws.Range("C53").Value = "Shepard Events"
$ws.Range("D53").Value = 3
$ws.Range("E53").Value = 1
$ws.Range("A54").Value = 41
$ws.Range("C54").Value = "Shepard Exposition Services"
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 0
$ws.Range("A55").Value = 36
$ws.Range("C55").Value = "Suburban Indy Home & Outdoor Living Shows"
$ws.Range("D55").Value = 3.8
$ws.Range("E55").Value = 13
$ws.Range("A56").Value = 49
$ws.Range("C56").Value = "The Indiana Convention center"
$ws.Range("D56").Value = 4.7
$ws.Range("E56").Value = 21
$ws.Range("A57").Value = 55
$ws.Range("C57").Value = "The Korner Garage"
$ws.Range("D57").Value = 4.5
$ws.Range("E57").Value = 6
$ws.Range("A58").Value = 42
$ws.Range("C58").Value = "West Pavilion"
$ws.Range("D58").Value = 4.3
$ws.Range("E58").Value = 33
$ws.Range("A59").Value = 53
$ws.Range("C59").Value = "Whale of a Sale"
$ws.Range("D59").Value = 4.7
$ws.Range("E59").Value = 26
